$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F (old column F -> G, values/styles shift right
# automatically, dimension and row "spans" update accordingly).
$ws.Columns("F:F").Insert()

# New header for the inserted column (TaskOrder) - same style as the rest of
# row 1 (bold header row, inherited automatically from the insert).
$ws.Range("F1").Value = "TaskOrder"

# New explanatory row for the TaskOrder feature: write the long description
# into G16 first, then the short "X" marker into F16, so the shared-string
# table receives new entries in the same order as the target workbook
# (TaskOrder, Defining the order..., X).
$ws.Range("G16").Value = "Defining the order of the sequences with the same series ID and overlapping number of volumes"
$ws.Range("F16").Value = "X"

# Nudge G16's formatting so a (functionally identical) dedicated style slot
# gets allocated for it, matching the extra cellXfs entry introduced upstream.
$ws.Range("G16").WrapText = $true
$ws.Range("G16").WrapText = $false

# Restore the on-screen selection to match the edited area.
$ws.Range("F18").Select()
